# Generate Report for Handback
# Renames the in-flight handback file from 1442afcb... to 29a8db48...
# (updating its generated xliff names/dates) and appends a brand new
# handback row for 457d907f-d578-46a9-a135-862072b6db8b.md on every sheet.

$wb = $excel.ActiveWorkbook

$renamedGuid = "29a8db48-3038-46fa-a4f9-36e163ec02d9"
$newGuid = "457d907f-d578-46a9-a135-862072b6db8b"

$renamedHash = "41322cac0cab4eb5d4073b0ed8fb1c0fd8bcc28c"
$newHash = "a7ed9e9cebfe4bd18a6c94ce7ec14bf903368b74"

$statusText = "Handed back: in sync with en-US"
$dateFmt = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$wsOv = $wb.Worksheets.Item("Overview")

# Row 2: file got renamed
$wsOv.Range("A2").Value2 = "$renamedGuid.md"
$wsOv.Range("C2").Value2 = ".md"
$wsOv.Range("E2").Value2 = $statusText
$wsOv.Range("F2").Value2 = $statusText
$wsOv.Range("G2").Value2 = "2016-09-03 19:14:51"
$wsOv.Range("G2").NumberFormat = $dateFmt

# Row 3: brand new file handed back
$wsOv.Range("A3").Value2 = "$newGuid.md"
$wsOv.Range("C3").Value2 = ".md"
$wsOv.Range("E3").Value2 = $statusText
$wsOv.Range("F3").Value2 = $statusText
$wsOv.Range("G3").Value2 = "2016-09-03 19:14:51"
$wsOv.Range("G3").NumberFormat = $dateFmt

# hyperlinks (B2 renamed, B3 new) -- rebuild all hyperlinks on this sheet
$wsOv.Hyperlinks.Delete()
$wsOv.Hyperlinks.Add($wsOv.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b1a6c9c1b0f0d9f4a1c9e9b3a9f5e9c9d9f9a9c1/e2e/$renamedGuid.md", "", "", "e2e\$renamedGuid.md")
$wsOv.Hyperlinks.Add($wsOv.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c2b7dad2c1e1eaf5b2dafacb0a6f0fad0eafbfc2/e2e/$newGuid.md", "", "", "e2e\$newGuid.md")

$loOv = $wsOv.ListObjects.Item(1)
$loOv.Resize($wsOv.Range("A1:G3"))

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

# Row 2: renamed file, regenerated xliff + dates
$wsZh.Range("B2").Value2 = ".md"
$wsZh.Range("C2").Value2 = $statusText
$wsZh.Range("D2").Value2 = "e2e"
$wsZh.Range("E2").Value2 = "ht"
$wsZh.Range("F2").Value2 = "'False"
$wsZh.Range("G2").Value2 = "$renamedGuid.$renamedHash.zh-cn.xlf"
$wsZh.Range("H2").Value2 = "2016-09-03 19:14:47"
$wsZh.Range("H2").NumberFormat = $dateFmt
$wsZh.Range("J2").Value2 = "$renamedGuid.$renamedHash.zh-cn.xlf"
$wsZh.Range("K2").Value2 = "2016-09-03 19:15:10"
$wsZh.Range("K2").NumberFormat = $dateFmt
$wsZh.Range("L2").Value2 = "'"
$wsZh.Range("M2").Value2 = "'True"
$wsZh.Range("N2").Value2 = "'"
$wsZh.Range("O2").Value2 = "'False"
$wsZh.Range("P2").Value2 = "'"

# Row 3: brand new file
$wsZh.Range("B3").Value2 = ".md"
$wsZh.Range("C3").Value2 = $statusText
$wsZh.Range("D3").Value2 = "e2e"
$wsZh.Range("E3").Value2 = "ht"
$wsZh.Range("F3").Value2 = "'True"
$wsZh.Range("G3").Value2 = "$newGuid.$newHash.zh-cn.xlf"
$wsZh.Range("H3").Value2 = "2016-09-03 19:14:47"
$wsZh.Range("H3").NumberFormat = $dateFmt
$wsZh.Range("J3").Value2 = "$newGuid.$newHash.zh-cn.xlf"
$wsZh.Range("K3").Value2 = "2016-09-03 19:15:10"
$wsZh.Range("K3").NumberFormat = $dateFmt
$wsZh.Range("L3").Value2 = "'"
$wsZh.Range("M3").Value2 = "'True"
$wsZh.Range("N3").Value2 = "'"
$wsZh.Range("O3").Value2 = "'False"
$wsZh.Range("P3").Value2 = "'"

# hyperlinks (A2,I2 renamed, A3,I3 new) -- rebuild all hyperlinks on this sheet
$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b1a6c9c1b0f0d9f4a1c9e9b3a9f5e9c9d9f9a9c1/e2e/$renamedGuid.md", "", "", "$renamedGuid.md")
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/d3c8eabe3d2fba06c3ebfbdc1b7a1fbe1fbcacd3/e2e/$renamedGuid.md", "", "", "$renamedGuid.md")
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c2b7dad2c1e1eaf5b2dafacb0a6f0fad0eafbfc2/e2e/$newGuid.md", "", "", "$newGuid.md")
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/e4d9fbcf4e3fcb17d4fcfcedfc8b2fcf2fcdbde4/e2e/$newGuid.md", "", "", "$newGuid.md")

$loZh = $wsZh.ListObjects.Item(1)
$loZh.Resize($wsZh.Range("A1:P3"))

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

# Row 2: renamed file, regenerated xliff + dates
$wsDe.Range("B2").Value2 = ".md"
$wsDe.Range("C2").Value2 = $statusText
$wsDe.Range("D2").Value2 = "e2e"
$wsDe.Range("E2").Value2 = "ht"
$wsDe.Range("F2").Value2 = "'False"
$wsDe.Range("G2").Value2 = "$renamedGuid.$renamedHash.de-de.xlf"
$wsDe.Range("H2").Value2 = "2016-09-03 19:14:51"
$wsDe.Range("H2").NumberFormat = $dateFmt
$wsDe.Range("J2").Value2 = "$renamedGuid.$renamedHash.de-de.xlf"
$wsDe.Range("K2").Value2 = "2016-09-03 19:15:18"
$wsDe.Range("K2").NumberFormat = $dateFmt
$wsDe.Range("L2").Value2 = "'"
$wsDe.Range("M2").Value2 = "'True"
$wsDe.Range("N2").Value2 = "'"
$wsDe.Range("O2").Value2 = "'False"
$wsDe.Range("P2").Value2 = "'"

# Row 3: brand new file
$wsDe.Range("B3").Value2 = ".md"
$wsDe.Range("C3").Value2 = $statusText
$wsDe.Range("D3").Value2 = "e2e"
$wsDe.Range("E3").Value2 = "ht"
$wsDe.Range("F3").Value2 = "'True"
$wsDe.Range("G3").Value2 = "$newGuid.$newHash.de-de.xlf"
$wsDe.Range("H3").Value2 = "2016-09-03 19:14:51"
$wsDe.Range("H3").NumberFormat = $dateFmt
$wsDe.Range("J3").Value2 = "$newGuid.$newHash.de-de.xlf"
$wsDe.Range("K3").Value2 = "2016-09-03 19:15:18"
$wsDe.Range("K3").NumberFormat = $dateFmt
$wsDe.Range("L3").Value2 = "'"
$wsDe.Range("M3").Value2 = "'True"
$wsDe.Range("N3").Value2 = "'"
$wsDe.Range("O3").Value2 = "'False"
$wsDe.Range("P3").Value2 = "'"

# hyperlinks (A2,I2 renamed, A3,I3 new) -- rebuild all hyperlinks on this sheet
$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b1a6c9c1b0f0d9f4a1c9e9b3a9f5e9c9d9f9a9c1/e2e/$renamedGuid.md", "", "", "$renamedGuid.md")
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/f5eafcdaf5f4cdc28e5fdcfcfd9c3fdf3fdecef5/e2e/$renamedGuid.md", "", "", "$renamedGuid.md")
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c2b7dad2c1e1eaf5b2dafacb0a6f0fad0eafbfc2/e2e/$newGuid.md", "", "", "$newGuid.md")
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/06fbfdebf06f5dedf9f6fedfded0f4fe04fefdf0/e2e/$newGuid.md", "", "", "$newGuid.md")

$loDe = $wsDe.ListObjects.Item(1)
$loDe.Resize($wsDe.Range("A1:P3"))
